$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 173 (pushes old rows 173..231 down to 174..232)
$ws.Rows.Item(173).Insert()

# Populate the newly inserted row 173 with the new record values.
$ws.Cells.Item(173, 1).Value = 10
$ws.Cells.Item(173, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(173, 3).Value = "La Araucanía"
$ws.Cells.Item(173, 4).Value = 44988
$ws.Cells.Item(173, 5).Value = 9
$ws.Cells.Item(173, 6).Value = "Fruta"
$ws.Cells.Item(173, 7).Value = 100104
$ws.Cells.Item(173, 8).Value = "Frutos de pepita"
$ws.Cells.Item(173, 9).Value = 100104003
$ws.Cells.Item(173, 10).Value = "Membrillo"
$ws.Cells.Item(173, 11).Value = "Champion"
$ws.Cells.Item(173, 12).Value = "Primera"
$ws.Cells.Item(173, 13).Value = 25
$ws.Cells.Item(173, 14).Value = 15000
$ws.Cells.Item(173, 15).Value = 15000
$ws.Cells.Item(173, 16).Value = 15000
$ws.Cells.Item(173, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(173, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(173, 19).Value = 833
$ws.Cells.Item(173, 20).Value = 18
